$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, [string]$text) {
    # Force the cell to stay text (so date-looking strings like "02-12-2025"
    # are not auto-converted to date serials), write the value, then reset
    # the style back to Normal so we don't leave a stray custom number
    # format behind.
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    if ($text -eq "") {
        # A bare "" clears the cell entirely in Excel. Use a lone
        # quote-prefix so the cell keeps its Text type with an empty value,
        # matching how the other blank cells in this sheet are stored.
        $r.Value = "'"
    } else {
        $r.Value = $text
    }
    $r.Style = "Normal"
}

# --- Row 7: the 02-12-2025 Medha Sub Division Office Coll. cash deposit was
# actually rolled into the Monday entry, so this row's own deposit date and
# amount are cleared out.
Set-TextCell "D7" ""
$ws.Range("E7").Value = 0

# --- Append 4 more collection rows (12-15) after the existing data.
$ws.Rows("12:15").Insert()

Set-TextCell "A12" "02-12-2025"
Set-TextCell "B12" "020965017-Kai Lalsingrao Shinde Gr.Bid.S.S.Pat.Ltd Kudal Br. Kudal"
Set-TextCell "C12" "Cash"
Set-TextCell "D12" "2025-12-02"
$ws.Range("E12").Value = 70200
Set-TextCell "F12" ""
Set-TextCell "G12" "2025-12-23"

Set-TextCell "A13" "02-12-2025"
Set-TextCell "B13" "020965019-SHRI DATTATRAY MAHARAJ KALAMBE SAH. PAT. LTD.DAPAWADI"
Set-TextCell "C13" "Cash"
Set-TextCell "D13" "2025-12-02"
$ws.Range("E13").Value = 17700
Set-TextCell "F13" ""
Set-TextCell "G13" "2025-12-23"

Set-TextCell "A14" "02-12-2025"
Set-TextCell "B14" "020965020-KAI.LALSINGRAO BAPUSO SHINDE SAH.PAT.LTD.,KUDAL BR.SAYGAON"
Set-TextCell "C14" "Cash"
Set-TextCell "D14" "2025-12-02"
$ws.Range("E14").Value = 13010
Set-TextCell "F14" ""
Set-TextCell "G14" "2025-12-23"

Set-TextCell "A15" "02-12-2025"
Set-TextCell "B15" "010965012-Medha Sub Division Office Coll."
Set-TextCell "C15" "Cash"
Set-TextCell "D15" ""
$ws.Range("E15").Value = 0
Set-TextCell "F15" "दिनांक 02.12.2025 रोजी रविवार असल्याने जमा झालेली रक्कम ही सोमवार दिनांक 03.12.2025 रोजी बँकेमध्ये भरणा करण्यात आली."
Set-TextCell "G15" "2025-12-23"
